# edit.ps1 - applies "grammar and edit to report" changes to hw_1_report.docx
#
# Strategy: Word's Find/Execute handles simple text edits fine, but several
# of the required edits need to *split* existing runs and interleave new
# <w:proofErr/> and <w:bookmarkStart/<w:bookmarkEnd/> markers between them -
# something the high level Range.Text setter can't express. For those
# paragraphs we clear the paragraph's run content (keeping its <w:pPr/>)
# and re-insert a literal OOXML fragment for the paragraph via
# Range.InsertXML, which Word's COM layer accepts as a raw WordML snippet.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-ParagraphXml {
    param($Paragraph, $PPrXml, $BodyXml)
    $rng = $Paragraph.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Delete()
    $xml = "<w:p $wNs>$PPrXml$BodyXml</w:p>"
    $rng.InsertXML($xml) | Out-Null
}

# ---------------------------------------------------------------------
# 1) Drop the stray leading " " run in the "Missing: (...)" bullet.
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -match '^\s?Missing: \(ACCTAGE') {
        $target = $d.Paragraphs($i)
        break
    }
}
$pr = $target.Range
$lead = $d.Range($pr.Start, $pr.Start + 1)
if ($lead.Text -eq " ") {
    $lead.Delete()
}

# ---------------------------------------------------------------------
# 2) "foring savings ... insurance." bullet: insert a gramStart/gramEnd
#    around "So" and split " ... checking account." into four runs
#    separated by " or checking" / " account".
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match '^It would be interesting to compare whether') {
        $target = $d.Paragraphs($i)
        break
    }
}

$pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$body = ''
$body += '<w:r><w:t>It would be interesting to compare whether individuals who have a checking or a savings account would b</w:t></w:r>'
$body += '<w:r><w:t xml:space="preserve">e more likely to get insurance.  The odds ratios </w:t></w:r>'
$body += '<w:proofErr w:type="spellStart"/>'
$body += '<w:r><w:t>foring</w:t></w:r>'
$body += '<w:proofErr w:type="spellEnd"/>'
$body += '<w:r><w:t xml:space="preserve"> savings and checking accounts are 1.779 and 0.379, respectively.  This would suggest that customers who have a savings account 1.779 more than those without a savings account to get insurance while customers with a checking account are only 0.379 more likely to get a checking account than those without a checking account.  </w:t></w:r>'
$body += '<w:proofErr w:type="gramStart"/>'
$body += '<w:r><w:t>So</w:t></w:r>'
$body += '<w:proofErr w:type="gramEnd"/>'
$body += '<w:r><w:t xml:space="preserve"> if all you knew about a potential customer if </w:t></w:r>'
$body += '<w:proofErr w:type="spellStart"/>'
$body += '<w:r><w:t>if</w:t></w:r>'
$body += '<w:proofErr w:type="spellEnd"/>'
$body += '<w:r><w:t xml:space="preserve"> they have a savings account</w:t></w:r>'
$body += '<w:r><w:t xml:space="preserve"> or checking</w:t></w:r>'
$body += '<w:r><w:t xml:space="preserve"> account</w:t></w:r>'
$body += '<w:r><w:t xml:space="preserve"> you could guess customers with savings account are more likely to get insurance.</w:t></w:r>'

Replace-ParagraphXml $target $pPr $body

# ---------------------------------------------------------------------
# 3) "In total 3034 observations had a missing value." bullet: merge the
#    two runs that used to sandwich the (now relocated) _GoBack bookmark.
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match '^In total 3034 observat') {
        $target = $d.Paragraphs($i)
        break
    }
}
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}

$pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$body = '<w:r><w:t>In total 3034 observations had a missing value.</w:t></w:r>'
Replace-ParagraphXml $target $pPr $body

# ---------------------------------------------------------------------
# 4) "We removed variables ... insignificant." bullet: append " as alpha
#    = 0.5" as its own run, then relocate the _GoBack bookmark here
#    (right before a trailing "." run).
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match '^We removed variables with missing values') {
        $target = $d.Paragraphs($i)
        break
    }
}

$pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$body = ''
$body += '<w:r><w:t>We removed variables with missing values rather than rows which kept the same total number of observations in our dataset, but information could have been lost in the dropped variables.  Looking at the dropped variables’ histograms there are some variables that would interesting to including</w:t></w:r>'
$body += '<w:r><w:t xml:space="preserve"> such as income, age, age of oldest account, and credit score.  Credit score in particular would be worth keeping in the model considering it only has 195 missing observations.  Testing this returns a model with an AIC of 9857.6 which is worse than our previous best model and the CRSCORE estimate is also insignificant</w:t></w:r>'
$body += '<w:r><w:t xml:space="preserve"> as alpha = 0.5</w:t></w:r>'
$body += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$body += '<w:r><w:t>.</w:t></w:r>'

Replace-ParagraphXml $target $pPr $body

# ---------------------------------------------------------------------
# 5) styles.xml: DefaultParagraphFont gains <w:semiHidden/>.
# ---------------------------------------------------------------------
# (handled separately - see below)
